$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-LatticeCell($row, $col, $innerXml) {
    $cell = $tbl.Cell($row, $col)
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' + $innerXml + '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$cell.Range.InsertXML($xml)
}

Set-LatticeCell 1 1 '<w:t>76 x 33</w:t><w:br/><w:t xml:space="preserve">  3    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>6|    |</w:t>'
Set-LatticeCell 1 2 '<w:t>89 x 22</w:t><w:br/><w:t xml:space="preserve">  2    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>9|    |</w:t>'
Set-LatticeCell 1 3 '<w:t>18 x 10</w:t><w:br/><w:t xml:space="preserve">  1    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>8|    |</w:t>'
Set-LatticeCell 2 1 '<w:t>71 x 48</w:t><w:br/><w:t xml:space="preserve">  4    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>1|    |</w:t>'
Set-LatticeCell 2 2 '<w:t>22 x 60</w:t><w:br/><w:t xml:space="preserve">  6    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>2|    |</w:t>'
Set-LatticeCell 2 3 '<w:t>40 x 27</w:t><w:br/><w:t xml:space="preserve">  2    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>0|    |</w:t>'
Set-LatticeCell 3 1 '<w:t>34 x 42</w:t><w:br/><w:t xml:space="preserve">  4    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>4|    |</w:t>'
Set-LatticeCell 3 2 '<w:t>70 x 85</w:t><w:br/><w:t xml:space="preserve">  8    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>0|    |</w:t>'
Set-LatticeCell 3 3 '<w:t>29 x 69</w:t><w:br/><w:t xml:space="preserve">  6    9</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>9|    |</w:t>'
Set-LatticeCell 4 1 '<w:t>10 x 78</w:t><w:br/><w:t xml:space="preserve">  7    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>0|    |</w:t>'
Set-LatticeCell 4 2 '<w:t>90 x 71</w:t><w:br/><w:t xml:space="preserve">  7    1</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>0|    |</w:t>'
Set-LatticeCell 4 3 '<w:t>65 x 19</w:t><w:br/><w:t xml:space="preserve">  1    9</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>5|    |</w:t>'
Set-LatticeCell 5 1 '<w:t>15 x 55</w:t><w:br/><w:t xml:space="preserve">  5    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>5|    |</w:t>'
Set-LatticeCell 5 2 '<w:t>65 x 24</w:t><w:br/><w:t xml:space="preserve">  2    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>5|    |</w:t>'
Set-LatticeCell 5 3 '<w:t>43 x 39</w:t><w:br/><w:t xml:space="preserve">  3    9</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>3|    |</w:t>'
